$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 122 (weekly update), shifting all subsequent rows down by one.
$ws.Rows.Item(122).Insert()

# Populate the new row 122 with this week's new price observation.
$ws.Cells.Item(122, 1).Value  = 10
$ws.Cells.Item(122, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(122, 3).Value  = "La Araucanía"
$ws.Cells.Item(122, 4).Value  = 44719
$ws.Cells.Item(122, 5).Value  = 9
$ws.Cells.Item(122, 6).Value  = 100112005
$ws.Cells.Item(122, 7).Value  = "Puerro"
$ws.Cells.Item(122, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(122, 9).Value  = "Primera"
$ws.Cells.Item(122, 10).Value = 40
$ws.Cells.Item(122, 11).Value = 12000
$ws.Cells.Item(122, 12).Value = 13000
$ws.Cells.Item(122, 13).Value = 12500
$ws.Cells.Item(122, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(122, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(122, 16).Value = 1042
$ws.Cells.Item(122, 17).Value = 12
$ws.Cells.Item(122, 18).Value = "Hortaliza"
